$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fixtures data for rows 2-14 (Fixture, Pick, Competition, Time, Win Info, Confidence %, Odds)
$data = @(
    @("Fulham vs Arsenal", "Arsenal", "England", "18th Oct 17:30", "81 of 106 Tips", "76", "1.52"),
    @("Torino vs Napoli", "Napoli", "Italy", "18th Oct 17:00", "71 of 78 Tips", "91", "1.85"),
    @("Marseille vs Le Havre", "Marseille", "France", "18th Oct 20:05", "51 of 54 Tips", "94", "1.38"),
    @("Atletico Madrid vs Osasuna", "Atletico Madrid", "Spain", "18th Oct 20:00", "50 of 54 Tips", "93", "1.38"),
    @("Angers vs Monaco", "Monaco", "France", "18th Oct 18:00", "37 of 42 Tips", "88", "1.60"),
    @("Astana vs Aktobe", "Astana", "Kazakhstan", "19th Oct 13:00", "32 of 33 Tips", "97", "1.41"),
    @("Bayern Munich vs Borussia Dortmund", "Bayern Munich", "Germany", "18th Oct 17:30", "25 of 30 Tips", "83", "1.40"),
    @("Liverpool vs Man Utd", "Liverpool", "England", "19th Oct 16:30", "24 of 34 Tips", "71", "1.61"),
    @("Cienciano vs Cusco FC", "Cienciano", "Peru", "19th Oct 00:00", "22 of 23 Tips", "96", "2.35"),
    @("Getafe vs Real Madrid", "Real Madrid", "Spain", "19th Oct 20:00", "20 of 23 Tips", "87", "1.46"),
    @("FC Ordabasy vs Okzhetpes", "FC Ordabasy", "Kazakhstan", "19th Oct 14:00", "18 of 18 Tips", "100", "1.73"),
    @("AC Milan vs Fiorentina", "AC Milan", "Italy", "19th Oct 19:45", "16 of 19 Tips", "84", "1.55"),
    @("Barcelona vs Olympiacos", "Barcelona", "Europe", "21st Oct 17:45", "15 of 18 Tips", "83", "1.18")
)

# Confidence %/Odds columns hold text-formatted numbers (matches source feed format);
# format as text BEFORE writing so the values are stored as strings, not numbers.
$ws.Range("F2:G14").NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $ws.Range("F$row").Value = $r[5]
    $ws.Range("G$row").Value = $r[6]
    $row++
}

# Remove the now-unused trailing rows (old rows 15-18)
$ws.Range("A15:H18").ClearContents()

# Re-place the summary formula, now one row up, referencing the shrunk range
$ws.Range("H15").Formula = "=AVERAGE(F2:F14)"
# Undo any formatting picked up from the adjacent text column
$ws.Range("H15").Style = "Normal"
